$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the shared style (bold font, thin box border, center/top alignment)
# on B1, then copy/paste the formatting onto A2 so both cells end up
# referencing the very same cell style (avoids creating extra unused styles).
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1   # xlContinuous
$b1.Borders.Weight = 2      # xlThin
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop

$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
